$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New user rows (8-12), mirroring the existing "User::factory()" table.
# ---------------------------------------------------------------------------

# Row 8 - Laralila Aguilar
$ws.Range("A8").Value = "laralila.aguilar@salud.qroo.gob.mx"
$ws.Range("B8").Value = "Laralila Aguilar"
$ws.Range("C8").Value = "laguilar"
$ws.Range("D8").Formula = '="\App\Models\User::factory()->create(["'
$ws.Range("E8").Formula = '="''name''=>''"&B8&"'',''email''=>''"&A8&"'',''password''=>''"&C8&"''"'
$ws.Range("F8").Formula = '=D8&E8&"]);"'

# Row 9 - Ricardo May
$ws.Range("A9").Value = "ricardo.may@salud.qroo.gob.mx"
$ws.Range("B9").Value = "Ricardo May"
$ws.Range("C9").Value = "rmay"
$ws.Range("D9").Formula = '="\App\Models\User::factory()->create(["'
$ws.Range("E9").Formula = '="''name''=>''"&B9&"'',''email''=>''"&A9&"'',''password''=>''"&C9&"''"'
$ws.Range("F9").Formula = '=D9&E9&"]);"'

# Row 10 - Margarita Che
$ws.Range("A10").Value = "margarita.che@salud.qroo.gob.mx"
$ws.Range("B10").Value = "Margarita Che"
$ws.Range("C10").Value = "mche"
$ws.Range("D10").Formula = '="\App\Models\User::factory()->create(["'
$ws.Range("E10").Formula = '="''name''=>''"&B10&"'',''email''=>''"&A10&"'',''password''=>''"&C10&"''"'
$ws.Range("F10").Formula = '=D10&E10&"]);"'

# Row 11 - Veneralda Rosado
$ws.Range("A11").Value = "veneralda.rosado@salud.qroo.gob.mx"
$ws.Range("B11").Value = "Veneralda Rosado"
$ws.Range("C11").Value = "vrosado"
$ws.Range("D11").Formula = '="\App\Models\User::factory()->create(["'
$ws.Range("E11").Formula = '="''name''=>''"&B11&"'',''email''=>''"&A11&"'',''password''=>''"&C11&"''"'
$ws.Range("F11").Formula = '=D11&E11&"]);"'

# Row 12 - Norman Angulo
$ws.Range("A12").Value = "norman.angulo@salud.qroo.gob.m"
$ws.Range("B12").Value = "Norman Angulo"
$ws.Range("C12").Value = "nangulo"
$ws.Range("D12").Formula = '="\App\Models\User::factory()->create(["'
$ws.Range("E12").Formula = '="''name''=>''"&B12&"'',''email''=>''"&A12&"'',''password''=>''"&C12&"''"'
$ws.Range("F12").Formula = '=D12&E12&"]);"'

# ---------------------------------------------------------------------------
# Formatting for the new rows.
# ---------------------------------------------------------------------------

# A8 gets a "plain" style: just the small grey Arial font used throughout
# the sheet, no fill/border/alignment changes.
$rngA8 = $ws.Range("A8")
$fA8 = $rngA8.Font
$fA8.Name = "Arial"
$fA8.Size = 8
$fA8.Color = 3355443

# B8 gets the "boxed" style used by all the B/C columns in the new rows:
# small grey Arial font, white fill, medium light-grey left/right border,
# centered + wrapped text.
$rngB8 = $ws.Range("B8")
$fB8 = $rngB8.Font
$fB8.Name = "Arial"
$fB8.Size = 8
$fB8.Color = 3355443
$rngB8.Interior.Pattern = 1
$rngB8.Interior.Color = 16777215
$rngB8.Borders.Item(7).LineStyle = 1
$rngB8.Borders.Item(7).Weight = -4138
$rngB8.Borders.Item(7).Color = 14540253
$rngB8.Borders.Item(10).LineStyle = 1
$rngB8.Borders.Item(10).Weight = -4138
$rngB8.Borders.Item(10).Color = 14540253
$rngB8.HorizontalAlignment = -4108
$rngB8.VerticalAlignment = -4108
$rngB8.WrapText = $true

# Re-use the B8 "boxed" style (copy formats only) on every other B/C cell
# in the new rows.
$ws.Range("B8").Copy()
$ws.Range("C8,B9,C9,B10,C10,B11,C11,B12,C12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# A10 gets the "boxed" style too, but with the full (4-sided) medium grey
# border already used elsewhere on the sheet, and vertical-only centering.
$rngA10 = $ws.Range("A10")
$fA10 = $rngA10.Font
$fA10.Name = "Arial"
$fA10.Size = 8
$fA10.Color = 3355443
$rngA10.Interior.Pattern = 1
$rngA10.Interior.Color = 16777215
$rngA10.Borders.Item(7).LineStyle = 1
$rngA10.Borders.Item(7).Weight = -4138
$rngA10.Borders.Item(7).Color = 14540253
$rngA10.Borders.Item(10).LineStyle = 1
$rngA10.Borders.Item(10).Weight = -4138
$rngA10.Borders.Item(10).Color = 14540253
$rngA10.Borders.Item(8).LineStyle = 1
$rngA10.Borders.Item(8).Weight = -4138
$rngA10.Borders.Item(8).Color = 14540253
$rngA10.Borders.Item(9).LineStyle = 1
$rngA10.Borders.Item(9).Weight = -4138
$rngA10.Borders.Item(9).Color = 14540253
$rngA10.VerticalAlignment = -4108
$rngA10.WrapText = $true

# ---------------------------------------------------------------------------
# Misc housekeeping to match the final saved state of the workbook.
# ---------------------------------------------------------------------------
$ws.Range("E6").Select()

Write-Host "done"
